$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 16
$ws.Range("D1").Value = $false

$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = $false

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = $false
